# Daily update at 8 AM UTC
# Appends the next day's row of results to the "Wins Over Time" tracking
# sheet, and moves the "most-recent row" date style (date-only, no time)
# from the old last row to the newly appended last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the current last used row in column A (currently row 46).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

# The number format used by every "interior" (non-last) row in column A.
$interiorFormat = $ws.Cells.Item(2, 1).NumberFormat
# The special number format used by the last row in column A.
$lastRowFormat = $ws.Cells.Item($lastRow, 1).NumberFormat

# Append the new day's data.
$ws.Cells.Item($newRow, 1).Value = 45787
$ws.Cells.Item($newRow, 2).Value = 192
$ws.Cells.Item($newRow, 3).Value = 203
$ws.Cells.Item($newRow, 4).Value = 197

# The old last row (A46) is no longer the last row, so it reverts to the
# regular interior date format ...
$ws.Cells.Item($lastRow, 1).NumberFormat = $interiorFormat
# ... and the new last row (A47) takes on the special last-row date format.
$ws.Cells.Item($newRow, 1).NumberFormat = $lastRowFormat
